# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (H-N) for a set of leve rows
# across several crafting job sheets, matching the latest market data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 408.5
$ws.Range("I31").Value = 408.5
$ws.Range("K31").Value = 1225.5
$ws.Range("M31").Value = -995.5

# Row 106
$ws.Range("H106").Value = 25211.688
$ws.Range("I106").Value = 25211.688
$ws.Range("K106").Value = 25211.688
$ws.Range("M106").Value = -24580.688

# Row 138
$ws.Range("H138").Value = 2067.0942
$ws.Range("I138").Value = 853.63635
$ws.Range("J138").Value = 2384.9048
$ws.Range("K138").Value = 2560.90905
$ws.Range("L138").Value = 7154.714399999999
$ws.Range("M138").Value = 2579.09095
$ws.Range("N138").Value = -17434.7144

$ws = $wb.Worksheets.Item("ARM")
# Row 9
$ws.Range("H9").Value = 20007.5
$ws.Range("J9").Value = 20007
$ws.Range("L9").Value = 20007
$ws.Range("N9").Value = -20347

# Row 20
$ws.Range("H20").Value = 20007.5
$ws.Range("J20").Value = 20007
$ws.Range("L20").Value = 20007
$ws.Range("N20").Value = -20547

# Row 32
$ws.Range("H32").Value = 1468.6522
$ws.Range("I32").Value = 1103.9318
$ws.Range("J32").Value = 9492.5
$ws.Range("K32").Value = 1103.9318
$ws.Range("L32").Value = 9492.5
$ws.Range("M32").Value = -816.9318000000001
$ws.Range("N32").Value = -10066.5

# Row 37
$ws.Range("H37").Value = 23430.715
$ws.Range("I37").Value = 25756.25
$ws.Range("J37").Value = 20330
$ws.Range("K37").Value = 25756.25
$ws.Range("L37").Value = 20330
$ws.Range("M37").Value = -25483.25
$ws.Range("N37").Value = -20876

# Row 110
$ws.Range("H110").Value = 1934
$ws.Range("I110").Value = 1914.875
$ws.Range("J110").Value = 2010.5
$ws.Range("K110").Value = 1914.875
$ws.Range("L110").Value = 2010.5
$ws.Range("M110").Value = 130.125
$ws.Range("N110").Value = -6100.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1010.1111
$ws.Range("I94").Value = 1053.75
$ws.Range("K94").Value = 1053.75
$ws.Range("M94").Value = -602.75

# Row 112
$ws.Range("H112").Value = 58888
$ws.Range("J112").Value = 58888
$ws.Range("L112").Value = 58888
$ws.Range("N112").Value = -61842

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3936.6428
$ws.Range("I58").Value = 4606.4
$ws.Range("J58").Value = 3564.5557
$ws.Range("K58").Value = 4606.4
$ws.Range("L58").Value = 3564.5557
$ws.Range("M58").Value = -4403.4
$ws.Range("N58").Value = -3970.5557

# Row 107
$ws.Range("H107").Value = 1133
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 1149.5
$ws.Range("K107").Value = 1100
$ws.Range("L107").Value = 1149.5
$ws.Range("M107").Value = 820
$ws.Range("N107").Value = -4989.5

# Row 136
$ws.Range("H136").Value = 3936.6428
$ws.Range("I136").Value = 4606.4
$ws.Range("J136").Value = 3564.5557
$ws.Range("K136").Value = 13819.2
$ws.Range("L136").Value = 10693.6671
$ws.Range("M136").Value = -11269.2
$ws.Range("N136").Value = -15793.6671

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 1919.875
$ws.Range("I13").Value = 1371.8
$ws.Range("K13").Value = 4115.4
$ws.Range("M13").Value = -3947.4

# Row 23
$ws.Range("H23").Value = 598.1667
$ws.Range("I23").Value = 104
$ws.Range("J23").Value = 845.25
$ws.Range("K23").Value = 312
$ws.Range("L23").Value = 2535.75
$ws.Range("M23").Value = -77
$ws.Range("N23").Value = -3005.75

# Row 54
$ws.Range("H54").Value = 7666
$ws.Range("J54").Value = 5999.5
$ws.Range("L54").Value = 17998.5
$ws.Range("N54").Value = -19116.5

# Row 62
$ws.Range("H62").Value = 17000
$ws.Range("J62").Value = 17000
$ws.Range("L62").Value = 51000
$ws.Range("N62").Value = -52372

# Row 65
$ws.Range("H65").Value = 17000
$ws.Range("J65").Value = 17000
$ws.Range("L65").Value = 153000
$ws.Range("N65").Value = -159864

# Row 68
$ws.Range("H68").Value = 1701.4
$ws.Range("I68").Value = 1384.6
$ws.Range("K68").Value = 4153.799999999999
$ws.Range("M68").Value = -3342.799999999999

# Row 71
$ws.Range("H71").Value = 1701.4
$ws.Range("I71").Value = 1384.6
$ws.Range("K71").Value = 12461.4
$ws.Range("M71").Value = -8405.4

# Row 80
$ws.Range("H80").Value = 5149.625
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 5699.5
$ws.Range("K80").Value = 10500
$ws.Range("L80").Value = 17098.5
$ws.Range("M80").Value = -9564
$ws.Range("N80").Value = -18970.5

# Row 83
$ws.Range("H83").Value = 5149.625
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 5699.5
$ws.Range("K83").Value = 31500
$ws.Range("L83").Value = 51295.5
$ws.Range("M83").Value = -26820
$ws.Range("N83").Value = -60655.5

# Row 122
$ws.Range("H122").Value = 853.0909
$ws.Range("I122").Value = 796.7273
$ws.Range("J122").Value = 909.4545000000001
$ws.Range("K122").Value = 7170.545700000001
$ws.Range("L122").Value = 8185.0905
$ws.Range("M122").Value = -4720.545700000001
$ws.Range("N122").Value = -13085.0905

# Row 124
$ws.Range("H124").Value = 1399
$ws.Range("I124").Value = 998.5
$ws.Range("J124").Value = 1666
$ws.Range("K124").Value = 2995.5
$ws.Range("L124").Value = 4998
$ws.Range("M124").Value = 1914.5
$ws.Range("N124").Value = -14818

# Row 140
$ws.Range("H140").Value = 1691.8667
$ws.Range("J140").Value = 2499.8
$ws.Range("L140").Value = 7499.400000000001
$ws.Range("N140").Value = -17859.4

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 10012125
$ws.Range("I11").Value = 25018876
$ws.Range("J11").Value = 2508749.5
$ws.Range("K11").Value = 25018876
$ws.Range("L11").Value = 2508749.5
$ws.Range("M11").Value = -25018737
$ws.Range("N11").Value = -2509027.5

# Row 12
$ws.Range("H12").Value = 250
$ws.Range("J12").Value = 250
$ws.Range("L12").Value = 250
$ws.Range("N12").Value = -530

# Row 24
$ws.Range("H24").Value = 10309.5
$ws.Range("J24").Value = 10102.615
$ws.Range("L24").Value = 10102.615
$ws.Range("N24").Value = -10448.615

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3213.375
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# Row 27
$ws.Range("H27").Value = 3213.375
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 20007
$ws.Range("J18").Value = 20007
$ws.Range("L18").Value = 20007
$ws.Range("N18").Value = -20353

# Row 20
$ws.Range("H20").Value = 20011
$ws.Range("J20").Value = 20011
$ws.Range("L20").Value = 20011
$ws.Range("N20").Value = -20491

# Row 122
$ws.Range("H122").Value = 15154931
$ws.Range("I122").Value = 16132313
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 48396939
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -48394489
$ws.Range("N122").Value = -21400

# Row 136
$ws.Range("H136").Value = 5663.522
$ws.Range("I136").Value = 4628.778
$ws.Range("J136").Value = 9388.6
$ws.Range("K136").Value = 13886.334
$ws.Range("L136").Value = 28165.8
$ws.Range("M136").Value = -11336.334
$ws.Range("N136").Value = -33265.8

